$wb = $excel.ActiveWorkbook

# --- Controllers sheet: two progress values moved from 5% to 100% ---
$wsControllers = $wb.Worksheets.Item("Controllers")
$wsControllers.Range("D84").Value = 1
$wsControllers.Range("D85").Value = 1

# --- Daos sheet: one progress value moved from 0% to 100% ---
$wsDaos = $wb.Worksheets.Item("Daos")
$wsDaos.Range("C90").Value = 1

# --- Vistas sheet: several progress values updated ---
$wsVistas = $wb.Worksheets.Item("Vistas")
$wsVistas.Range("C8").Value = 1
$wsVistas.Range("C73").Value = 0.5
$wsVistas.Range("C74").Value = 1
$wsVistas.Range("C75").Value = 1

# --- Avance sheet: PEP estimate updated ---
$wsAvance = $wb.Worksheets.Item("Avance")
$wsAvance.Range("D11").Value = 14.5

# Recalculate so dependent formulas (COUNTIFS, SUM, shared formulas, etc.) update
$excel.Calculate()

# --- Update each sheet's view/selection state ---
$wsControllers.Activate()
$wsControllers.Range("D86").Select()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1

$wsDaos.Activate()
$wsDaos.Range("C94").Select()
$excel.ActiveWindow.ScrollRow = 85
$excel.ActiveWindow.ScrollColumn = 1

# Vistas becomes the active/selected tab (workbook activeTab points at it)
$wsVistas.Activate()
$wsVistas.Range("C74").Select()
$excel.ActiveWindow.ScrollRow = 63
$excel.ActiveWindow.ScrollColumn = 1
